$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 87 (existing rows 87-89 shift down to 88-90)
$ws.Rows.Item(87).Insert()

# Copy the number format (date) of the old D87 cell (now D88) into the new D87 cell
$ws.Range("D88").Copy()
$ws.Range("D87").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row 87 with the new record's values
$ws.Cells.Item(87, 1).Value = 3
$ws.Cells.Item(87, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44509
$ws.Cells.Item(87, 5).Value = 5
$ws.Cells.Item(87, 6).Value = 100112026
$ws.Cells.Item(87, 7).Value = "Haba"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 85
$ws.Cells.Item(87, 11).Value = 8000
$ws.Cells.Item(87, 12).Value = 8500
$ws.Cells.Item(87, 13).Value = 8235
$ws.Cells.Item(87, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(87, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(87, 16).Value = 329
$ws.Cells.Item(87, 17).Value = 25
$ws.Cells.Item(87, 18).Value = "Hortaliza"

$wb.Save()
